# Restore previously-removed time-interval rows (Wenglor instrument
# readings) to the bottom of the InstrumentMetadata_Jericoacoara sheet,
# and update the sheet view's selection to point at the newly active area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New data rows (107-114). Columns:
# A Site, B Date, C StartTime, D EndTime, E InstrumentType, F Instrument,
# G StartHeight_m, H EndHeight_m, I HeightErr_m, J HeightRef,
# K Longitudinal_m, L Spanwise_m, M AngleErr_deg, N ErrorCode, O InstrumentID
# ---------------------------------------------------------------------
$rows = @(
    @{ R=107; A="Jericoacoara"; B=41963; C="0.472222222222222";  D="0.513888888888889"; E="Wenglor"; F="W3"; G="-0.155";                   H="-0.155";                   I="0.001"; J="L2"; K="0"; L="-0.079";                  M="0"; N="1"; O="21" },
    @{ R=108; A="Jericoacoara"; B=41963; C="0.472222222222222";  D="0.513888888888889"; E="Wenglor"; F="W4"; G="-0.0040000000000000036";   H="-0.0040000000000000036";   I="0.001"; J="L2"; K="0"; L="-0.079";                  M="0"; N="0"; O="14" },
    @{ R=109; A="Jericoacoara"; B=41963; C="0.472222222222222";  D="0.513888888888889"; E="Wenglor"; F="W5"; G="-0.097";                   H="-0.097";                   I="0.001"; J="S2"; K="0"; L="-0.496";                  M="0"; N="0"; O="R6" },
    @{ R=110; A="Jericoacoara"; B=41963; C="0.472222222222222";  D="0.513888888888889"; E="Wenglor"; F="W6"; G="-0.097";                   H="-0.097";                   I="0.001"; J="S2"; K="0"; L="-0.665";                  M="0"; N="0"; O="R4" },
    @{ R=111; A="Jericoacoara"; B=41963; C="0.472222222222222";  D="0.513888888888889"; E="Wenglor"; F="W7"; G="-0.097";                   H="-0.097";                   I="0.001"; J="S2"; K="0"; L="-0.929";                  M="0"; N="0"; O="R7" },
    @{ R=112; A="Jericoacoara"; B=41963; C="0.472222222222222";  D="0.513888888888889"; E="Wenglor"; F="W8"; G="0.08099999999999999";      H="0.08099999999999999";      I="0.001"; J="S2"; K="0"; L="-0.929";                  M="0"; N="0"; O="R5" },
    @{ R=113; A="Jericoacoara"; B=41963; C="0.47222222222222227"; D="0.513888888888889"; E="Wenglor"; F="W1"; G="-0.155";                   H="-0.155";                   I="0.001"; J="L2"; K="0"; L="0.251";                   M="0"; N="1"; O="R9" },
    @{ R=114; A="Jericoacoara"; B=41963; C="0.47222222222222227"; D="0.513888888888889"; E="Wenglor"; F="W2"; G="-0.155";                   H="-0.155";                   I="0.001"; J="L2"; K="0"; L="0.082";                   M="0"; N="1"; O="R8" }
)

# Instrument IDs that look numeric ("21", "14") get a quote-prefixed text
# style in the source data (as opposed to the alphanumeric codes like R6-R9)
$quotedIds = @("21", "14")

foreach ($row in $rows) {
    $r = $row.R

    $ws.Range("A$r").Value = $row.A

    $ws.Range("B$r").NumberFormat = "yyyy\-mm\-dd"
    $ws.Range("B$r").Value = $row.B

    $ws.Range("C$r").NumberFormat = "hh:mm:ss.00"
    $ws.Range("C$r").Value = [double]$row.C

    $ws.Range("D$r").NumberFormat = "hh:mm:ss.00"
    $ws.Range("D$r").Value = [double]$row.D

    $ws.Range("E$r").NumberFormat = "hh:mm:ss.00"
    $ws.Range("E$r").Value = $row.E

    $ws.Range("F$r").Value = $row.F

    $ws.Range("G$r").Value = [double]$row.G
    $ws.Range("H$r").Value = [double]$row.H
    $ws.Range("I$r").Value = [double]$row.I

    $ws.Range("J$r").Value = $row.J

    $ws.Range("K$r").Value = [double]$row.K
    $ws.Range("L$r").Value = [double]$row.L
    $ws.Range("M$r").Value = [double]$row.M
    $ws.Range("N$r").Value = [double]$row.N

    $ws.Range("O$r").NumberFormat = "@"
    if ($quotedIds -contains $row.O) {
        $ws.Range("O$r").Formula = "'" + $row.O
    } else {
        $ws.Range("O$r").Value = $row.O
    }
}

# ---------------------------------------------------------------------
# Update the frozen-pane view: scroll the bottom pane so row 86 is the
# top-most visible row, and move the active selection to E93.
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 86
$win.ScrollColumn = 1
$ws.Range("E93").Select()
